$wb = $excel.ActiveWorkbook

# --- Update status text and timestamps (Generate Report for Handoff) ---

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-26 09:04:58"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-26 09:04:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-26 09:04:58"

# --- Narrow the "status" columns that used to be sized for the long
#     "Handed back: in sync with en-US" text, now that the text is shorter ---

$narrowWidth = 16.3333333333333

$wsOverview.Range("E1").ColumnWidth = $narrowWidth
$wsOverview.Range("F1").ColumnWidth = $narrowWidth

$wsZhCn.Range("C1").ColumnWidth = $narrowWidth

$wsDeDe.Range("C1").ColumnWidth = $narrowWidth
